$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.68%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.34%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.795"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.92%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06255"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.51%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.839"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.56%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8773"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.10%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9459"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.64%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1464"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.40%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05168"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.99%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07267"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.47%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03151"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.77%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.04%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001553"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.24%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006256"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.59%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005759"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.34%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.477"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.72%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.261"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.79%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.275"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.09%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.63%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1309"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.10%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.843"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.10%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04323"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.85%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001174"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.82%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004257"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.33%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001198"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.23%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001685"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.17%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04024"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.04%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006238"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "51.67%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1150"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.50%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002123"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.14%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01352"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.88%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005134"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.27%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.22%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.930"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "2,870.93%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02987"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-12.24%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002097"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.22%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001997"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.22%"
